$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (VC): update Measured value + add Error formula ---
$ws.Range("D3").Value = 7.6
$ws.Range("E3").Formula = "=ABS(C3-D3)/C3"

# --- Row 4 (VB): update Measured value ---
$ws.Range("D4").Value = 2.072

# --- Row 5 (VE): update Measured value ---
$ws.Range("D5").Value = 1.465

# --- Row 6 (VCE): D6 is an existing shared formula (=D3-D5); it recalculates
#     automatically now that D3 & D5 changed, so it is left untouched.

# --- Row 7 (Ib): D7 becomes a formula now ---
$ws.Range("D7").Formula = "=14.12*10^-6"

# --- Row 8 (Ic): update Measured value ---
$ws.Range("D8").Value = 0.00287

# --- Row 9: D9 formula (SUM) unchanged, recalculates automatically ---

# --- New column E (Error) shared formula across rows 4-11 (only 4-9 keep cells) ---
$ws.Range("E4:E11").Formula = "=ABS(C4-D4)/C4"
$ws.Range("E10").ClearContents()
$ws.Range("E11").ClearContents()

# --- Rows 10 & 11: clear the IR1 / IR2 rows (keep C/D styles, drop labels & values) ---
$ws.Range("A10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()

# --- Row 13: Max Gain measured value changes ---
$ws.Range("D13").Value = 152

# --- Row 14: F-3dB Dominant measured value becomes text "5.1k" ---
$ws.Range("D14").Value = "5.1k"

# --- Row 15: Fmax measured value text changes from 9.7k to 9.9k ---
$ws.Range("D15").Value = "9.9k"

# --- Row 16: F-3dB 2 measured value text changes from 375K to 480k ---
$ws.Range("D16").Value = "480k"

# --- Row 1: new "Error" header in column E ---
$ws.Range("E1").Value = "Error"

# --- Row 19 (new row): Frequency Bandwidth theoretical/measured ---
$ws.Range("A19").Value = "Frequency Bandwidth"
$ws.Range("C19").Formula = "=20971-1.22"
$ws.Range("D19").Formula = "=480-5.1"

# --- Update selection / scroll position to match the saved view ---
$ws.Range("E21").Select()
